# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) in selected bullet points,
# matching the commit: "Implement quantitative metrics highlighting
# across all resume formats".
#
# Strategy: for each target paragraph, walk left-to-right through the
# paragraph's text range (excluding the trailing paragraph mark) and use
# Find.Execute, constrained to a shrinking sub-range, to locate each
# numeric token in order. Each found token's Range gets Font.Bold = true
# and Font.Color set to the corporate-blue accent (RGB 2C3E50, encoded
# as a BGR COM color value). This causes Word to split the run exactly
# around the matched text, which reproduces the target OOXML run
# structure (plain-text runs interleaved with bold+colored runs).

$d = $word.ActiveDocument

# RGB 2C3E50 expressed as a Word/COM BGR color long (0x00502E2C... actually
# 0x00503E2C) -> wdColor value 5258796 renders as w:color w:val="2C3E50".
$accentColor = 5258796

function BoldToken($doc, $rangeStart, $rangeEnd, $needle) {
    $searchRange = $doc.Range($rangeStart, $rangeEnd)
    $found = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "token not found: $needle"
    }
    $searchRange.Font.Bold = $true
    $searchRange.Font.Color = $accentColor
    return $searchRange.End
}

function HighlightParagraph($doc, $paragraphIndex, $tokens) {
    $para = $doc.Paragraphs.Item($paragraphIndex)
    $full = $para.Range
    $paraStart = $full.Start
    $paraEnd = $full.End
    # Exclude the paragraph mark itself from the searchable limit so that
    # no formatting ever gets written into the paragraph-mark run (w:pPr/w:rPr).
    $limit = $paraEnd - 1

    $cursor = $paraStart
    foreach ($tok in $tokens) {
        $cursor = BoldToken $doc $cursor $limit $tok
    }
}

# 1) Partner - Siege Analytics bullet: demographic classification accuracy
HighlightParagraph $d 10 @('23%', '64%')

# 2) Partner - Siege Analytics bullet: sampling / survey margin of error
HighlightParagraph $d 12 @('±4.2%', '±2.1%', '71%', '87%')

# 3) Partner - Siege Analytics bullet: boundary estimation mapping costs
HighlightParagraph $d 13 @('73.5%', '$4.7M')

# 4) Partner - Siege Analytics bullet: FEC analysis / political spending
HighlightParagraph $d 14 @('$2')

# 5) Data Products Manager bullet: ETL processing time reduction
HighlightParagraph $d 24 @('57%')

# 6) Key Achievements: revenue generation
HighlightParagraph $d 50 @('$4.9M')

# 7) Key Achievements: conversion rate improvement
HighlightParagraph $d 51 @('23%')

# 8) Key Achievements: platform impact analyst count
HighlightParagraph $d 53 @('12,847')

Write-Output "Highlighting applied to 8 paragraphs."
